$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column (H) mirrors the existing header styling (bold, centered,
# bordered) used by the other header cells (e.g. G1 "sum"), so copy that
# cell's formatting over before writing the new header text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 1
